$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$green = 5287936

# --- Clean slate: unmerge + clear the whole area we will rebuild ---
$ws.Range("A5:L24").UnMerge()
$ws.Range("A5:L24").Clear()

# --- Cell values ---
$ws.Range("B5").Value = "При нормальных условиях"
$ws.Range("A6").Value = "Размер презентации"
$ws.Range("B6").Value = "2 Мб, 12 слайдов"
$ws.Range("C6").Value = "30 Мб, 31 слайд"
$ws.Range("D6").Value = "300 слайдов"
$ws.Range("E6").Value = "1000 слайдов"
$ws.Range("A7").Value = "CPU(%)"
$ws.Range("B7").Value = "27"
$ws.Range("C7").Value = "32"
$ws.Range("D7").Value = "32"
$ws.Range("E7").Value = "32"
$ws.Range("A8").Value = "Памати(Мбайт)"
$ws.Range("B8").Value = "190"
$ws.Range("C8").Value = "190"
$ws.Range("D8").Value = "280"
$ws.Range("E8").Value = "280"
$ws.Range("B9").Value = "При 1 свободном процессоре"
$ws.Range("A10").Value = "Состояние"
$ws.Range("B10").Value = "корректно"
$ws.Range("C10").Value = "корректно"
$ws.Range("D10").Value = "корректно"
$ws.Range("A11").Value = "Комментарий"
$ws.Range("D11").Value = "В первые 10 секунд были подвисания, power point не отвечал"
$ws.Range("F11").Value = "Так как нагрузка при увеличении размера презентации не меняется, не имеет смысла тестировать на презентациях более 300 слайдов"
$ws.Range("B12").Value = "При 256 Мб оперативной памяти"
$ws.Range("F12").Value = "Один загруженный процессор это 25% от 100% ресурсов проыессоров"
$ws.Range("A13").Value = "Состояние"
$ws.Range("B13").Value = "корректно"
$ws.Range("C13").Value = "корректно"
$ws.Range("D13").Value = "корректно"
$ws.Range("F13").Value = "При меньших ресурсах не имеет смысла тестировать, так как не хватит ресуров даже на открытия окна Publish"
$ws.Range("A14").Value = "Комментарий"
$ws.Range("B14").Value = "Около 3 секунд были подвисания, power point не отвечал"
$ws.Range("D14").Value = "Около 5 секунд были подвисания, power point не отвечал"
$ws.Range("B15").Value = "При 1 свободном процессоре и 256 Мб оперативной памяти"
$ws.Range("A16").Value = "Состояние"
$ws.Range("B16").Value = "корректно"
$ws.Range("C16").Value = "корректно"
$ws.Range("D16").Value = "корректно"
$ws.Range("A17").Value = "Комментарий"
$ws.Range("B17").Value = "Около 6 секунд были подвисания, power point не отвечал"
$ws.Range("D17").Value = "Около 12 секунд были подвисания, power point не отвечал"

# --- Merges ---
$ws.Range("B5:E5").Merge()
$ws.Range("B9:E9").Merge()
$ws.Range("F13:L14").Merge()
$ws.Range("B12:E12").Merge()
$ws.Range("D14:E14").Merge()
$ws.Range("B14:C14").Merge()
$ws.Range("B22:E22").Merge()
$ws.Range("F11:L11").Merge()
$ws.Range("D10:E10").Merge()
$ws.Range("D11:E11").Merge()
$ws.Range("F12:L12").Merge()
$ws.Range("B15:E15").Merge()
$ws.Range("B17:C17").Merge()
$ws.Range("D17:E17").Merge()

# --- Styles (grouped by style signature, applied per-cell via loop) ---
foreach ($addr in @("A7","B7","C7","A8","B8","C8","D8","E8")) {
  $r = $ws.Range($addr)
}

foreach ($addr in @("F7","G7","H7")) {
  $r = $ws.Range($addr)
  $r.HorizontalAlignment = $xlCenter
  $r.VerticalAlignment = $xlCenter
  $r.WrapText = $true
}

foreach ($addr in @("L7")) {
  $r = $ws.Range($addr)
  $r.WrapText = $true
}

foreach ($addr in @("D7","E7")) {
  $r = $ws.Range($addr)
}

foreach ($addr in @("A5","A6","A9","A10","A11","A12","A13","A14","A15","A16","A17","A23","A24")) {
  $r = $ws.Range($addr)
  $r.Font.Bold = $true
}

foreach ($addr in @("B6","C6","D6","E6")) {
  $r = $ws.Range($addr)
  $r.Font.Bold = $true
  $r.HorizontalAlignment = $xlCenter
  $r.VerticalAlignment = $xlCenter
}

foreach ($addr in @("B10","C10","B13","C13","B16","C16")) {
  $r = $ws.Range($addr)
  $r.Interior.Color = $green
}

foreach ($addr in @("D13","E13","D16","E16")) {
  $r = $ws.Range($addr)
  $r.Interior.Color = $green
  $r.WrapText = $true
}

foreach ($addr in @("C12","D12","E12","F14","G14","H14","I14","J14","K14","L14")) {
  $r = $ws.Range($addr)
  $r.HorizontalAlignment = $xlCenter
  $r.VerticalAlignment = $xlCenter
  $r.WrapText = $true
}

foreach ($addr in @("D10","E10")) {
  $r = $ws.Range($addr)
  $r.Interior.Color = $green
  $r.WrapText = $true
}

foreach ($addr in @("F11","G11","H11","I11","J11","K11","L11","F12","G12","H12","I12","J12","K12","L12","F13","G13","H13","I13","J13","K13","L13","B14","C14","D14","E14","B17","C17","D17","E17")) {
  $r = $ws.Range($addr)
  $r.Font.Italic = $true
  $r.HorizontalAlignment = $xlCenter
  $r.VerticalAlignment = $xlCenter
  $r.WrapText = $true
}

foreach ($addr in @("B5","B9","B12","B15","B22")) {
  $r = $ws.Range($addr)
  $r.Font.Bold = $true
  $r.HorizontalAlignment = $xlCenter
  $r.VerticalAlignment = $xlCenter
  $r.WrapText = $true
}

foreach ($addr in @("C5","D5","E5","C9","D9","E9","C22","D22","E22")) {
  $r = $ws.Range($addr)
  $r.Font.Bold = $true
  $r.WrapText = $true
}

foreach ($addr in @("D11","E11")) {
  $r = $ws.Range($addr)
  $r.Font.Italic = $true
  $r.HorizontalAlignment = $xlCenter
  $r.WrapText = $true
}

foreach ($addr in @("C15","D15","E15")) {
  $r = $ws.Range($addr)
  $r.WrapText = $true
}

# --- Row heights ---
$ws.Rows.Item(7).RowHeight = 21
$ws.Rows.Item(9).RowHeight = 42.75
$ws.Rows.Item(11).RowHeight = 49.5
$ws.Rows.Item(12).RowHeight = 40.5
$ws.Rows.Item(13).RowHeight = 15
$ws.Rows.Item(14).RowHeight = 59.25
$ws.Rows.Item(15).RowHeight = 42
$ws.Rows.Item(17).RowHeight = 45.75
$ws.Rows.Item(18).RowHeight = 15
$ws.Rows.Item(20).RowHeight = 33.75

# --- Selection ---
$ws.Range("D14:E14").Select()
